$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text-formatted numbers (e.g. "41.789.78", "0.403").
# Force the whole Price column to Text so Excel does not silently convert
# values that happen to look like plain numbers into real numeric cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '41.805.97'
$ws.Range('E2').Value = '  +5.42%  '
$ws.Range('D3').Value = '2.228.49'
$ws.Range('E3').Value = '  +2.61%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '231.91'
$ws.Range('E5').Value = '  +2.22%  '
$ws.Range('D6').Value = '0.627'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '61.57'
$ws.Range('E7').Value = '  -2.38%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.403'
$ws.Range('E9').Value = '  +2.76%  '
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('E11').Value = '  +5.90%  '
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').Value = '2.559.95'
$ws.Range('E13').Value = '  +2.68%  '
$ws.Range('D14').Value = '15.65'
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('D15').Value = '22.01'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '0.802'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '5.58'
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').Value = '2.243.46'
$ws.Range('E18').Value = '  +3.38%  '
$ws.Range('D19').Value = '41.692.72'
$ws.Range('E19').Value = '  +5.20%  '
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('D21').Value = '72.11'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('D23').Value = '250.02'
$ws.Range('E23').Value = '  +9.03%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.37'
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').Value = '2.39'
$ws.Range('E26').Value = '  +2.15%  '
$ws.Range('D27').Value = '9.67'
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('D29').Value = '167.22'
$ws.Range('E29').Value = '  -2.15%  '
$ws.Range('D30').Value = '19.99'
$ws.Range('E30').Value = '  +0.70%  '
$ws.Range('D31').Value = '1.41'
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('D32').Value = '2.66'
$ws.Range('E32').Value = '  -0.82%  '
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('D34').Value = '4.99'
$ws.Range('E34').Value = '  +6.19%  '
$ws.Range('D35').Value = '4.67'
$ws.Range('E35').Value = '  +3.32%  '
$ws.Range('E36').Value = '  +3.17%  '
$ws.Range('E37').Value = '  -4.50%  '
$ws.Range('E38').Value = '  -5.19%  '
$ws.Range('E39').Value = '  -0.97%  '
$ws.Range('D40').Value = '0.000256'
$ws.Range('E40').Value = '  +30.50%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('B42').Value = 'FTXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D42').Value = '4.87'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0239'
$ws.Range('E43').Value = '  +4.25%  '
$ws.Range('D44').Value = '8.59'
$ws.Range('E44').Value = '  +8.75%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = '1.23'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '0.0976'
$ws.Range('E46').Value = '  +6.26%  '
$ws.Range('E47').Value = '  -3.54%  '
$ws.Range('D48').Value = '1.481.43'
$ws.Range('E48').Value = '  -2.13%  '
$ws.Range('D49').Value = '2.82'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('D50').Value = '16.52'
$ws.Range('D51').Value = '52.61'
$ws.Range('E51').Value = '  +7.06%  '

# Restore the default (unformatted) look now that the text values are stored,
# so no residual explicit cell formatting is left behind.
$priceRange.ClearFormats()
